$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target column-C (곡명/song title) values for rows 2..52, reflecting the
# re-scraped Naver Music Top 50 chart (one new entry, one dropped, reordered).
$songTitles = @(
    '곡명',
    'Dynamite',
    'Lovesick Girls',
    'DON''T TOUCH ME',
    '잠이 오질 않네요',
    '취기를 빌려 (취향저격 그녀 X 산들)',
    '힘든 건 사랑이 아니다',
    'I CAN''T STOP ME',
    'Savage Love (Laxed - Siren Beat) (BTS Remix)',
    '딩가딩가 (Dingga)',
    '마리아 (Maria)',
    'When We Disco (Duet with 선미)',
    '에잇(Prod.&Feat. SUGA of BTS)',
    '어떻게 이별까지 사랑하겠어, 널 사랑하는 거지',
    '눈누난나 (NUNU NANA)',
    'ALIEN',
    'Dolphin',
    '흔들리는 꽃들 속에서 네 샴푸향이 느껴진거야',
    '내 마음이 움찔했던 순간 (취향저격 그녀 X 규현)',
    'How You Like That',
    '놓아줘 (with 태연)',
    'Ice Cream (with Selena Gomez)',
    'Not Shy',
    '홀로',
    'Love poem',
    '너도 아는',
    '나랑 같이 걸을래 (바른연애 길잡이 X 적재)',
    '늦은 밤 너의 집 앞 골목길에서',
    '다시 여기 바닷가',
    '덤디덤디 (DUMDi DUMDi)',
    '아로하',
    'Blueming',
    '살짝 설렜어 (Nonstop)',
    'Make A Wish (Birthday Song)',
    '오래된 노래',
    'Bet You Wanna (Feat. Cardi B)',
    'Bad Boy',
    '도망가',
    '잠이 들어야 (Feat. 헤이즈)',
    '거짓말이라도 해서 널 보고싶어',
    '밤편지',
    'Memories',
    'Pretty Savage',
    '사랑하게 될 줄 알았어',
    '너를 만나',
    '보라빛 밤 (pporappippam)',
    'METEOR',
    'Downtown Baby',
    '모든 날, 모든 순간 (Every day, Every Moment)',
    'Holy (Feat. Chance The Rapper)',
    '우리 사랑한 동안'
)

for ($i = 0; $i -lt $songTitles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $songTitles[$i]
}

# New 51st chart row (row 52): rank index 50 in column A, sequence 51 in
# column B; column C already received its title above via the loop.
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 51

# Match A52 formatting (bold + border + centered/top) to the rest of column A.
$ws.Cells.Item(51, 1).Copy() | Out-Null
$ws.Cells.Item(52, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

